$d = $word.ActiveDocument

$replacements = @(
    @{old = "The system as built could have been a little better with a better follow-up (change of system responsible person)";
      new = "Design: The system as built could have been a little better with a better follow-up (change of system responsible person)"},
    @{old = "There was no rules clearly defining the need and requirements for this system.";
      new = "Design: There was no rules clearly defining the need and requirements for this system."},
    @{old = "System was purchased from MINIMAX who has been supplying the Clean agent system in all 3 shios of Vega series.";
      new = "Design: System was purchased from MINIMAX who has been supplying the Clean agent system in all 3 shios of Vega series."},
    @{old = "Some comments related to the arrangement of the room where equipment is located came on a late stage from Owner.";
      new = "Design: Some comments related to the arrangement of the room where equipment is located came on a late stage from Owner."},
    @{old = "Change of system responsible designer during the project made so that, there was a lack of communication and follow-up during the detail design phase.";
      new = "Design: Change of system responsible designer during the project made so that, there was a lack of communication and follow-up during the detail design phase."},
    @{old = "System was a little different due to change of supplier. Better arrangement around the system equipment was provided. ";
      new = "Design: System was a little different due to change of supplier. Better arrangement around the system equipment was provided. "}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
